# Apply Oct 12 2023 GitHub Actions cryptos-list price refresh to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.767.55'
$ws.Range('E2').Value = '  +0.02%  '
$ws.Range('D3').Value = '1.537.25'
$ws.Range('E3').Value = '  -1.89%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '205.36'
$ws.Range('E5').Value = '  -0.50%  '
$ws.Range('E6').Value = '  -0.80%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.246'
$ws.Range('E8').Value = '  -0.70%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '21.22'
$ws.Range('E9').Value = '  -3.05%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0580'
$ws.Range('E10').Value = '  -0.59%  '
$ws.Range('E11').Value = '  -0.91%  '
$ws.Range('D12').Value = '1.755.73'
$ws.Range('E12').Value = '  -1.85%  '
$ws.Range('D13').Value = '1.536.30'
$ws.Range('E13').Value = '  -2.04%  '
$ws.Range('E14').Value = '  -1.61%  '
$ws.Range('E15').Value = '  -1.47%  '
$ws.Range('D16').Value = '26.763.33'
$ws.Range('E16').Value = '  -0.08%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '60.88'
$ws.Range('E17').Value = '  -1.06%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '212.93'
$ws.Range('E18').Value = '  -0.81%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.22'
$ws.Range('E19').Value = '  -1.89%  '
$ws.Range('E20').Value = '  +0.46%  '
$ws.Range('E21').Value = '  -0.03%  '
$ws.Range('E22').Value = '  -2.04%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.14'
$ws.Range('E23').Value = '  -1.89%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.93'
$ws.Range('E24').Value = '  -3.11%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '151.86'
$ws.Range('E25').Value = '  -0.52%  '
$ws.Range('E26').Value = '  -2.51%  '
$ws.Range('E27').Value = '  -0.88%  '
$ws.Range('E28').Value = '  +0.04%  '
$ws.Range('E29').Value = '  -0.96%  '
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.10'
$ws.Range('E30').Value = '  -1.27%  '
$ws.Range('B31').Value = 'Hedera'
$ws.Range('C31').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0457'
$ws.Range('E31').Value = '  -1.44%  '
$ws.Range('E32').Value = '  +2.16%  '
$ws.Range('D33').Value = '1.363.50'
$ws.Range('E33').Value = '  -1.78%  '
$ws.Range('E34').Value = '  -0.23%  '
$ws.Range('E35').Value = '  -2.81%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.962'
$ws.Range('E36').Value = '  +4.22%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.28'
$ws.Range('E37').Value = '  -0.06%  '
$ws.Range('E38').Value = '  +0.68%  '
$ws.Range('E39').Value = '  -1.13%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.75'
$ws.Range('E40').Value = '  +7.71%  '
$ws.Range('E41').Value = '  -2.03%  '
$ws.Range('E43').Value = '  +0.76%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '62.86'
$ws.Range('E44').Value = '  -0.66%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.73'
$ws.Range('E45').Value = '  -3.59%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.26'
$ws.Range('E46').Value = '  -3.71%  '
$ws.Range('D47').Value = '1.670.35'
$ws.Range('E47').Value = '  -1.77%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '85.14'
$ws.Range('E48').Value = '  -0.49%  '
$ws.Range('E49').Value = '  +3.34%  '
$ws.Range('D50').Value = '0.0₇0978'
$ws.Range('E50').Value = '  -1.47%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0942'
$ws.Range('E51').Value = '  -1.04%  '
